$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Add the new "correlated_assigns_for_request (GET)" table at rows 24-25.
#    Set these cell values BEFORE renaming the row-17 title below, so that the
#    new shared strings get appended to the string table in the same order
#    the target workbook uses.
# ---------------------------------------------------------------------------
$ws.Range("A24").Value = "correlated_assigns_for_request (GET)"
$ws.Range("A25").Value = "id"
$ws.Range("B25").Value = "foreign_id(extractors.id) (for variable_name)"
$ws.Range("C25").Value = "foreign_id(requests.id)"
$ws.Range("D25").Value = "param_name"

# ---------------------------------------------------------------------------
# 2) Rename the existing "correlated_requests" table (row 17) to its new,
#    more descriptive title.
# ---------------------------------------------------------------------------
$ws.Range("A17").Value = "correlated_assigns_for_request (POST, PUT, fat GET, whatever)"

# ---------------------------------------------------------------------------
# 3) Apply formatting to the newly-added cells by copying it from the
#    equivalent, already-formatted cells of the neighboring table.
# ---------------------------------------------------------------------------
# Bold "table title" style (row 24 title, like A1 / A10 / A17)
$ws.Range("A17").Copy()
$ws.Range("A24").PasteSpecial(-4122) | Out-Null

# Bordered header-row style (row 25 headers, like row 18)
$ws.Range("A18:D18").Copy()
$ws.Range("A25:D25").PasteSpecial(-4122) | Out-Null

# Bordered empty data rows (rows 26-29), matching rows 19-22:
# row 19/26 has no value in column A, rows 20-22/27-29 do.
$ws.Range("B19:D19").Copy()
$ws.Range("B26:D26").PasteSpecial(-4122) | Out-Null

$ws.Range("A20:D20").Copy()
$ws.Range("A27:D27").PasteSpecial(-4122) | Out-Null

$ws.Range("A21:D21").Copy()
$ws.Range("A28:D28").PasteSpecial(-4122) | Out-Null

$ws.Range("A22:D22").Copy()
$ws.Range("A29:D29").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4) Update the window/selection state to match the edited workbook
#    (scrolled down a bit, with D19 as the active cell).
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 10
$null = $ws.Range("D19").Select()
